$d = $word.ActiveDocument

# 1) "ArcherProjectile - логика пули лучника" -> "ProjectilePool - создаем пул пуль "
$d.Content.Find.Execute("ArcherProjectile - логика пули лучника", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ProjectilePool - создаем пул пуль ", 2)

# 2) "Barracks - логика казармы" -> "ArcherPool - создаем пул лучников"
$d.Content.Find.Execute("Barracks - логика казармы", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ArcherPool - создаем пул лучников", 2)

# 3) Insert a brand-new bulleted/numbered paragraph right after the paragraph that
#    used to read "Barracks - логика казармы" (now "ArcherPool - создаем пул лучников"),
#    inheriting its numbering + run formatting, with the text
#    "ArcherProjectile - логика пуль лучников"
$found = $d.Content
$found.Find.Execute("ArcherPool - создаем пул лучников", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$archerPoolParagraph = $found.Paragraphs(1)
$archerPoolParagraph.Range.InsertParagraphAfter()
$newParagraph = $archerPoolParagraph.Next()
$newParagraph.Range.Text = "ArcherProjectile - логика пуль лучников"
